# Add new columns I (I0) and J (IF) to the sheet, mirroring the diff:
# - Header row: I1 = "I0", J1 = "IF" (same style as the existing header cells)
# - Data rows 2..80: numeric values for I and J (I and J hold identical values per row)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells - set text then clone the existing header's formatting (style index),
# instead of setting Font/Border/Alignment individually (which would create a new,
# slightly different style entry in styles.xml).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for column I (rows 2-80); column J mirrors column I for every row
$values = @(3,7,7,8,7,7,7,9,6,9,5,8,7,6,5,7,6,7,7,9,8,8,9,9,9,9,8,7,8,6,9,9,9,9,8,6,8,6,10,6,8,6,7,7,8,8,8,9,6,9,7,9,9,8,9,8,9,8,8,6,8,8,7,9,8,6,8,8,8,9,6,6,7,7,4,4,4,3,3)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $val = $values[$i]
    $ws.Cells.Item($row, 9).Value = $val   # Column I
    $ws.Cells.Item($row, 10).Value = $val  # Column J
}
